# Added the J and Z lines.
# Insert two new rows for the "Nassau Street" line (J/Z trains at Broad St
# and Bowery) just above the existing "Broadway" (49th St) block, shifting
# all subsequent rows (old 103-115) down to (105-117).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank rows at row 103 - this pushes the old rows 103-115
#    down to 105-117, preserving their contents.
$ws.Rows("103:104").Insert()

# 2) Copy the formatting (styles / row height) from an existing row that
#    has the same shape (exactly two populated line-code cells, D and E)
#    into the two freshly inserted blank rows so the new rows look like
#    every other data row in the table.
$ws.Range("A101:O101").Copy()
$ws.Range("A103:O103").PasteSpecial(-4122)
$ws.Range("A101:O101").Copy()
$ws.Range("A104:O104").PasteSpecial(-4122)
$ws.Rows("103:104").RowHeight = 20.35

# 3) Populate the two new rows with the Nassau Street (J/Z) station data.
$ws.Range("A103").Value = "Nassau Street"
$ws.Range("B103").Value = 103
$ws.Range("C103").Value = "Broad St"
$ws.Range("D103").Value = "J"
$ws.Range("E103").Value = "Z"
$ws.Range("F103").ClearContents()
$ws.Range("G103").ClearContents()

$ws.Range("A104").Value = "Nassau Street"
$ws.Range("B104").Value = 104
$ws.Range("C104").Value = "Bowery"
$ws.Range("D104").Value = "J"
$ws.Range("E104").Value = "Z"
$ws.Range("F104").ClearContents()
$ws.Range("G104").ClearContents()

# 4) Renumber column B ("stop #") for the rows that got shifted down, so
#    the sequence stays contiguous 105..117 instead of retaining the old
#    103..115 numbering.
for ($r = 105; $r -le 117; $r++) {
    $ws.Cells.Item($r, 2).Value = $r
}

# 5) Restore the view state: the active cell moves to C105 and the sheet
#    is scrolled down so row 102 is the first visible row.
$ws.Range("C105").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 102 } catch {}
try { $excel.ActiveWindow.TabRatio = 0.877 } catch {}
$null
